$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for even_MAG-GUT22049.fa (spreadsheet row 13), shifting subsequent rows up.
$ws.Rows(13).Delete()
